$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SHIP_ROUTES")
$ws.Activate()

# Insert a new column before column J (shifts J:P -> K:Q), duplicating column I's
# formatting, matching how Excel's "Insert Copied/Duplicate Column" works when the
# author duplicated Route 8 into a new Route 9.
$ws.Columns("J").Insert()

# New column J (Route 9) is a copy of old column I (Route 8) with a handful of
# differences: new header/ID, and Destination-3 fields that differ from Route 8.
$ws.Range("J1").Value = "Route 9"
$ws.Range("J2").Value = "North"
$ws.Range("J3").Value = 1.9
$ws.Range("J4").Value = "Gladstone"
$ws.Range("J5").Value = "FA"
$ws.Range("J6").Value = "FA_EXPSILO_STORE"
$ws.Range("J9").Value = "Melbourne"
$ws.Range("J10").Value = "FA"
$ws.Range("J11").Value = "FA_STORE"
$ws.Range("J14").Value = "Devonport"
$ws.Range("J15").Value = "GP"
$ws.Range("J16").Value = "GP_STORE"
$ws.Range("J17").Value = "Melbourne"
$ws.Range("J18").Value = "GP"
$ws.Range("J19").Value = "GP_STORE"
$ws.Range("J20").Value = "Gladstone"

# Column I's Destination-3 fields (rows 17-19) change from the old Route 8 values
# to new, distinct values now that Route 8 no longer shares column J.
$ws.Range("I17").Value = "Osborne"
$ws.Range("I18").Value = "FA"
$ws.Range("I19").Value = "FA_STORE"

$ws.Range("I7").Select()
